$d = $word.ActiveDocument

# --- 1. Replace the "Setback and buildable area for earthworks" picture
#        with a plain hyperlink run pointing at the image URL. ---
$url1 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C07_Earthworks.jpg?h=100%25&w=100%25"
$shp1 = $d.InlineShapes.Item(1)
$rng1 = $shp1.Range
$rng1.Text = $url1
$null = $d.Hyperlinks.Add($rng1, $url1)

# --- 2. Replace the "Dimensions and terracing of retaining walls" picture
#        with a plain hyperlink run pointing at the image URL. ---
$url2 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F15_Retaining_Wall_15m.jpg?h=100%2525&w=100%2525"
$shp2 = $d.InlineShapes.Item(1)
$rng2 = $shp2.Range
$rng2.Text = $url2
$null = $d.Hyperlinks.Add($rng2, $url2)

Write-Host "Remaining InlineShapes:" $d.InlineShapes.Count
